# device_configs.xlsx - rebuild sheet1 with connection-profile table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Wipe the old content (old A1:C3 "Baudrate / Config" table) completely so
#    we can lay out the new table from a clean sheet.
# ---------------------------------------------------------------------------
$ws.Cells.Clear() | Out-Null
$ws.Rows("2:3").Delete() | Out-Null

# ---------------------------------------------------------------------------
# 2. Header row (row 1)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "ConnectionType"
$ws.Range("B1").Value = "IP"
$ws.Range("C1").Value = "Port"
$ws.Range("D1").Value = "Username"
$ws.Range("E1").Value = "Password"
$ws.Range("F1").Value = "COMPort"
$ws.Range("G1").Value = "Baudrate"
$ws.Range("H1").Value = "DataBits"
$ws.Range("I1").Value = "Parity"
$ws.Range("J1").Value = "StopBits"
$ws.Range("K1").Value = "Config"

# ---------------------------------------------------------------------------
# 3. Row 2 - "serial" profile
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "serial"
$ws.Range("F2").Value = "COM1"
$ws.Range("G2").Value = 9600
$ws.Range("H2").Value = 8
$ws.Range("I2").Value = "None"
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = "enable`nconfiguration terminal"
$ws.Range("K2").WrapText = $true
$ws.Range("L2").WrapText = $true

# ---------------------------------------------------------------------------
# 4. Row 3 - "telnet" profile
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "telnet"
$ws.Range("B3").Value = "192.168.0.100"
$ws.Range("C3").Value = 23
$ws.Range("K3").Value = "ls"
$ws.Range("K3").WrapText = $true

# ---------------------------------------------------------------------------
# 5. Row 4 - "ssh" profile
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "ssh"
$ws.Range("B4").Value = "192.168.0.1"
$ws.Range("C4").Value = 22
$ws.Range("D4").Value = "admin"
$ws.Range("E4").Value = "admin"
$ws.Range("K4").Value = "cd .."
$ws.Range("K4").WrapText = $true

# ---------------------------------------------------------------------------
# 6. Row 6 - stray formatted (empty) cell
# ---------------------------------------------------------------------------
$ws.Range("K6").WrapText = $true

# ---------------------------------------------------------------------------
# 7. Number formats for the numeric-ish columns (Port, Baudrate, DataBits,
#    StopBits use a custom "0_);[Red](0)" format; COMPort column gets an
#    explicitly-applied General format)
# ---------------------------------------------------------------------------
$numFmt = "0_);[Red](0)"
$ws.Range("C1").NumberFormat = $numFmt
$ws.Range("C3").NumberFormat = $numFmt
$ws.Range("C4").NumberFormat = $numFmt
$ws.Range("G1").NumberFormat = $numFmt
$ws.Range("G2").NumberFormat = $numFmt
$ws.Range("H1").NumberFormat = $numFmt
$ws.Range("H2").NumberFormat = $numFmt
$ws.Range("J1").NumberFormat = $numFmt
$ws.Range("J2").NumberFormat = $numFmt

$ws.Range("F1").NumberFormat = "General"
$ws.Range("F2").NumberFormat = "General"

# ---------------------------------------------------------------------------
# 8. Column widths (final widths: A=16, B=12.875, C:J=10.625, K=21.75, L=21).
#    The inputs below are pre-compensated for the host's internal
#    char-width -> pixel -> char-width round-trip (+5/7) so the persisted
#    <col width="..."> ends up as close as possible to the intended value.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 15.285714285714286
$ws.Columns.Item(2).ColumnWidth = 12.142857142857142
$ws.Columns.Item(3).ColumnWidth = 9.857142857142858
$ws.Columns.Item(4).ColumnWidth = 9.857142857142858
$ws.Columns.Item(5).ColumnWidth = 9.857142857142858
$ws.Columns.Item(6).ColumnWidth = 9.857142857142858
$ws.Columns.Item(7).ColumnWidth = 9.857142857142858
$ws.Columns.Item(8).ColumnWidth = 9.857142857142858
$ws.Columns.Item(9).ColumnWidth = 9.857142857142858
$ws.Columns.Item(10).ColumnWidth = 9.857142857142858
$ws.Columns.Item(11).ColumnWidth = 21.0
$ws.Columns.Item(12).ColumnWidth = 20.285714285714285

# ---------------------------------------------------------------------------
# 9. Selection / active cell, matching the saved view state of the workbook
# ---------------------------------------------------------------------------
$ws.Range("F4").Select() | Out-Null
